$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 30 continues the "control d'affectation de AR (backend)" entry from
# row 29: same "algorithmique pur" solution, new remark, next day's date.

# Clone the border/font/number formats of row 29 (columns A-G) down into row 30.
$ws.Range("A29:G29").Copy()
$ws.Range("A30:G30").PasteSpecial(-4122)

# Column H on row 29 carries the hyperlink style (s=16); row 30's H has no
# hyperlink, so borrow the plain bordered style from G29 (s=14) instead.
$ws.Range("G29").Copy()
$ws.Range("H30").PasteSpecial(-4122)

$ws.Rows.Item(30).RowHeight = 114

$ws.Cells.Item(30, 2).Value = "shift to right if add shift to left if delete"
$ws.Cells.Item(30, 4).Value = "algorithmique pur"
$ws.Cells.Item(30, 5).Value = 42094

$ws.Range("A5:H30").Select()
